$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-07 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("605÷4=151, 1", $true, $false, $false, $false, $false, $true, 1, $false, "293÷8=36, 5", 2) | Out-Null
$d.Content.Find.Execute("178÷3=59, 1", $true, $false, $false, $false, $false, $true, 1, $false, "133÷2=66, 1", 2) | Out-Null
$d.Content.Find.Execute("759÷7=108, 3", $true, $false, $false, $false, $false, $true, 1, $false, "566÷9=62, 8", 2) | Out-Null
$d.Content.Find.Execute("751÷9=83, 4", $true, $false, $false, $false, $false, $true, 1, $false, "811÷3=270, 1", 2) | Out-Null
$d.Content.Find.Execute("198÷2=99, 0", $true, $false, $false, $false, $false, $true, 1, $false, "295÷7=42, 1", 2) | Out-Null
$d.Content.Find.Execute("234÷9=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "361÷6=60, 1", 2) | Out-Null
$d.Content.Find.Execute("220÷6=36, 4", $true, $false, $false, $false, $false, $true, 1, $false, "102÷9=11, 3", 2) | Out-Null
$d.Content.Find.Execute("355÷4=88, 3", $true, $false, $false, $false, $false, $true, 1, $false, "269÷4=67, 1", 2) | Out-Null
$d.Content.Find.Execute("842÷2=421, 0", $true, $false, $false, $false, $false, $true, 1, $false, "266÷6=44, 2", 2) | Out-Null
$d.Content.Find.Execute("514÷2=257, 0", $true, $false, $false, $false, $false, $true, 1, $false, "773÷9=85, 8", 2) | Out-Null
$d.Content.Find.Execute("728÷9=80, 8", $true, $false, $false, $false, $false, $true, 1, $false, "626÷2=313, 0", 2) | Out-Null
$d.Content.Find.Execute("248÷8=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "355÷3=118, 1", 2) | Out-Null
$d.Content.Find.Execute("986÷2=493, 0", $true, $false, $false, $false, $false, $true, 1, $false, "138÷9=15, 3", 2) | Out-Null
$d.Content.Find.Execute("250÷2=125, 0", $true, $false, $false, $false, $false, $true, 1, $false, "112÷5=22, 2", 2) | Out-Null
$d.Content.Find.Execute("768÷8=96, 0", $true, $false, $false, $false, $false, $true, 1, $false, "985÷3=328, 1", 2) | Out-Null
$d.Content.Find.Execute("115÷7=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "802÷7=114, 4", 2) | Out-Null
$d.Content.Find.Execute("129÷2=64, 1", $true, $false, $false, $false, $false, $true, 1, $false, "661÷2=330, 1", 2) | Out-Null
$d.Content.Find.Execute("321÷2=160, 1", $true, $false, $false, $false, $false, $true, 1, $false, "548÷3=182, 2", 2) | Out-Null
$d.Content.Find.Execute("398÷2=199, 0", $true, $false, $false, $false, $false, $true, 1, $false, "446÷8=55, 6", 2) | Out-Null
$d.Content.Find.Execute("569÷4=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "659÷6=109, 5", 2) | Out-Null
$d.Content.Find.Execute("588÷2=294, 0", $true, $false, $false, $false, $false, $true, 1, $false, "369÷6=61, 3", 2) | Out-Null
$d.Content.Find.Execute("758÷2=379, 0", $true, $false, $false, $false, $false, $true, 1, $false, "849÷4=212, 1", 2) | Out-Null
$d.Content.Find.Execute("772÷6=128, 4", $true, $false, $false, $false, $false, $true, 1, $false, "108÷3=36, 0", 2) | Out-Null
$d.Content.Find.Execute("969÷3=323, 0", $true, $false, $false, $false, $false, $true, 1, $false, "687÷8=85, 7", 2) | Out-Null
$d.Content.Find.Execute("174÷7=24, 6", $true, $false, $false, $false, $false, $true, 1, $false, "418÷8=52, 2", 2) | Out-Null
